# Update "想去人数" (interest/want-to-go counts) on the 展览, 演出, and 全部类型
# sheets to reflect the latest scraped values (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibition) sheet
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value = 1355
$wsExhibition.Range("F10").Value = 8420
$wsExhibition.Range("F11").Value = 456
$wsExhibition.Range("F19").Value = 10605
$wsExhibition.Range("F31").Value = 93
$wsExhibition.Range("F34").Value = 38

# 演出 (Performance) sheet
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F6").Value = 187

# 全部类型 (All Types) sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 187
$wsAll.Range("F5").Value = 1355
$wsAll.Range("F11").Value = 8420
$wsAll.Range("F12").Value = 456
$wsAll.Range("F20").Value = 10605
$wsAll.Range("F33").Value = 38
